$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data per commit diff

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '20.406.17'
$ws.Cells.Item(2, 5).Value = '  +2.51%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.462.26'
$ws.Cells.Item(3, 5).Value = '  +3.90%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.49%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '0.9499'
$ws.Cells.Item(5, 5).Value = '  -5.24%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '274.39'
$ws.Cells.Item(6, 5).Value = '  -0.48%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.3651'
$ws.Cells.Item(7, 5).Value = '  -0.15%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3076'
$ws.Cells.Item(8, 5).Value = '  -0.62%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -0.03%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '1.035'
$ws.Cells.Item(10, 5).Value = '  +0.44%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +0.94%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.9996'
$ws.Cells.Item(12, 5).Value = '  -0.38%  '

# Row 13
$ws.Cells.Item(13, 2).Value = 'Polkadot'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '5.415'
$ws.Cells.Item(13, 5).Value = '  -1.04%  '

# Row 14
$ws.Cells.Item(14, 2).Value = 'Solana'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '17.96'
$ws.Cells.Item(14, 5).Value = '  +2.36%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '6.117'
$ws.Cells.Item(15, 5).Value = '  -0.97%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.00001024'
$ws.Cells.Item(16, 5).Value = '  +0.71%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '1.462.33'
$ws.Cells.Item(17, 5).Value = '  +3.77%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.9663'
$ws.Cells.Item(18, 5).Value = '  -3.56%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.05778'
$ws.Cells.Item(19, 5).Value = '  +1.96%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '69.67'
$ws.Cells.Item(20, 5).Value = '  -1.55%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '5.421'
$ws.Cells.Item(21, 5).Value = '  -3.30%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '14.44'
$ws.Cells.Item(22, 5).Value = '  -1.74%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '10.86'
$ws.Cells.Item(23, 5).Value = '  -0.18%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  +0.23%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '20.438.39'
$ws.Cells.Item(25, 5).Value = '  +2.57%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '141.39'
$ws.Cells.Item(26, 5).Value = '  +6.59%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -7.59%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '17.12'
$ws.Cells.Item(28, 5).Value = '  -1.00%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.613.56'
$ws.Cells.Item(29, 5).Value = '  +2.81%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '112.08'
$ws.Cells.Item(30, 5).Value = '  +2.43%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '3.831'
$ws.Cells.Item(31, 5).Value = '  -1.87%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.866'
$ws.Cells.Item(32, 5).Value = '  -7.36%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'ImmutableX'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.7885'
$ws.Cells.Item(33, 5).Value = '  -3.16%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'Stellar'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.07804'
$ws.Cells.Item(34, 5).Value = '  +1.59%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +2.09%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.05688'
$ws.Cells.Item(36, 5).Value = '  -1.88%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '4.659'
$ws.Cells.Item(37, 5).Value = '  -4.91%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +3.66%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.02024'
$ws.Cells.Item(39, 5).Value = '  -1.67%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.9554'
$ws.Cells.Item(40, 5).Value = '  -4.56%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '10.34'
$ws.Cells.Item(41, 5).Value = '  -0.94%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'FraxShare'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '7.448'
$ws.Cells.Item(42, 5).Value = '  -10.26%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'Algorand'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.1856'
$ws.Cells.Item(43, 5).Value = '  -1.68%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.5257'
$ws.Cells.Item(44, 5).Value = '  -0.78%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -1.47%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '11.86'
$ws.Cells.Item(46, 5).Value = '  -4.06%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '116.94'
$ws.Cells.Item(47, 5).Value = '  +2.22%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.5136'
$ws.Cells.Item(48, 5).Value = '  -0.70%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.746'
$ws.Cells.Item(49, 5).Value = '  -0.96%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.06424'
$ws.Cells.Item(50, 5).Value = '  +4.23%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.9842'
$ws.Cells.Item(51, 5).Value = '  -1.85%  '
